# This script applies a full-row data reshuffle to the "Artfynd" sheet.
# Source data (columns A,B,D,E,F,G,H,Q,R and occasionally AC) for each
# row was rearranged among the rows, matching an upstream re-sync of the
# underlying species-observation export. Rows 6 and 23 are unaffected.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111486571
$ws.Range("B2").Value = 96348
$ws.Range("D2").Value = "VU"
$ws.Range("E2").Value = 220787
$ws.Range("F2").Value = "Knärot"
$ws.Range("G2").Value = "Goodyera repens"
$ws.Range("H2").Value = "(L.) R. Br."
$ws.Range("Q2").Value = 610186.5997174035
$ws.Range("R2").Value = 6897342.822581144
# Row 3
$ws.Range("A3").Value = 111486578
$ws.Range("B3").Value = 77267
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 6446
$ws.Range("F3").Value = "Kolflarnlav"
$ws.Range("G3").Value = "Carbonicola anthracophila"
$ws.Range("H3").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("Q3").Value = 610248.2770640558
$ws.Range("R3").Value = 6897273.826739896
# Row 4
$ws.Range("A4").Value = 111486564
$ws.Range("B4").Value = 5135
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 105930
$ws.Range("F4").Value = "Vågbandad barkbock"
$ws.Range("G4").Value = "Semanotus undatus"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
$ws.Range("Q4").Value = 610237.2946205279
$ws.Range("R4").Value = 6897509.394074276
# Row 5
$ws.Range("A5").Value = 111486559
$ws.Range("B5").Value = 96348
$ws.Range("D5").Value = "VU"
$ws.Range("E5").Value = 220787
$ws.Range("F5").Value = "Knärot"
$ws.Range("G5").Value = "Goodyera repens"
$ws.Range("H5").Value = "(L.) R. Br."
$ws.Range("Q5").Value = 610236.04261225
$ws.Range("R5").Value = 6897547.660441305
# Row 7
$ws.Range("A7").Value = 111486563
$ws.Range("B7").Value = 77186
$ws.Range("E7").Value = 353
$ws.Range("F7").Value = "Dvärgbägarlav"
$ws.Range("G7").Value = "Cladonia parasitica"
$ws.Range("H7").Value = "(Hoffm.) Hoffm."
$ws.Range("Q7").Value = 610237.2946205279
$ws.Range("R7").Value = 6897509.394074276
# Row 8
$ws.Range("A8").Value = 111486565
$ws.Range("B8").Value = 78081
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 229821
$ws.Range("F8").Value = "Vedflamlav"
$ws.Range("G8").Value = "Ramboldia elabens"
$ws.Range("H8").Value = "(Fr.) Kantvilas & Elix"
$ws.Range("Q8").Value = 610237.2946205279
$ws.Range("R8").Value = 6897509.394074276
# Row 9
$ws.Range("A9").Value = 111486566
$ws.Range("Q9").Value = 610196.0688674429
$ws.Range("R9").Value = 6897453.853924472
$ws.Range("AC9").Value = "på mossbeväxt block"
# Row 10
$ws.Range("A10").Value = 111486573
$ws.Range("B10").Value = 96383
$ws.Range("D10").Value = "LC"
$ws.Range("E10").Value = 223621
$ws.Range("F10").Value = "Skogsnattviol"
$ws.Range("G10").Value = "Platanthera bifolia subsp. latiflora"
$ws.Range("H10").Value = "(Drejer) Løjtnant"
$ws.Range("Q10").Value = 610146.8202491006
$ws.Range("R10").Value = 6897400.387088978
# Row 11
$ws.Range("A11").Value = 111486567
$ws.Range("B11").Value = 93057
$ws.Range("E11").Value = 2809
$ws.Range("F11").Value = "Mörk husmossa"
$ws.Range("G11").Value = "Hylocomiastrum umbratum"
$ws.Range("H11").Value = "(Hedw.) M.Fleisch."
$ws.Range("Q11").Value = 610195.0273142112
$ws.Range("R11").Value = 6897457.090065848
# Row 12
$ws.Range("A12").Value = 111486562
$ws.Range("B12").Value = 96348
$ws.Range("D12").Value = "VU"
$ws.Range("E12").Value = 220787
$ws.Range("F12").Value = "Knärot"
$ws.Range("G12").Value = "Goodyera repens"
$ws.Range("H12").Value = "(L.) R. Br."
$ws.Range("Q12").Value = 610236.2225055038
$ws.Range("R12").Value = 6897513.563481026
# Row 13
$ws.Range("A13").Value = 111486582
$ws.Range("B13").Value = 88489
$ws.Range("E13").Value = 1962
$ws.Range("F13").Value = "Vaddporing"
$ws.Range("G13").Value = "Anomoporia kamtschatica"
$ws.Range("H13").Value = "(Parmasto) Bondartseva"
$ws.Range("Q13").Value = 610142.0730319817
$ws.Range("R13").Value = 6897316.605782338
# Row 14
$ws.Range("A14").Value = 111486586
$ws.Range("B14").Value = 96348
$ws.Range("D14").Value = "VU"
$ws.Range("E14").Value = 220787
$ws.Range("F14").Value = "Knärot"
$ws.Range("G14").Value = "Goodyera repens"
$ws.Range("H14").Value = "(L.) R. Br."
$ws.Range("Q14").Value = 610319.2657305499
$ws.Range("R14").Value = 6897606.443173738
# Row 15
$ws.Range("A15").Value = 111486570
$ws.Range("Q15").Value = 610138.220419018
$ws.Range("R15").Value = 6897377.214110524
# Row 16
$ws.Range("A16").Value = 111486572
$ws.Range("B16").Value = 96348
$ws.Range("D16").Value = "VU"
$ws.Range("E16").Value = 220787
$ws.Range("F16").Value = "Knärot"
$ws.Range("G16").Value = "Goodyera repens"
$ws.Range("H16").Value = "(L.) R. Br."
$ws.Range("Q16").Value = 610186.866762756
$ws.Range("R16").Value = 6897391.885662847
# Row 17
$ws.Range("A17").Value = 111486552
$ws.Range("B17").Value = 78081
$ws.Range("E17").Value = 229821
$ws.Range("F17").Value = "Vedflamlav"
$ws.Range("G17").Value = "Ramboldia elabens"
$ws.Range("H17").Value = "(Fr.) Kantvilas & Elix"
$ws.Range("Q17").Value = 610202.9845822605
$ws.Range("R17").Value = 6897614.31648521
# Row 18
$ws.Range("A18").Value = 111486585
$ws.Range("B18").Value = 73634
$ws.Range("E18").Value = 6426
$ws.Range("F18").Value = "Kattfotslav"
$ws.Range("G18").Value = "Felipes leucopellaeus"
$ws.Range("H18").Value = "(Ach.) Frisch & G.Thor"
$ws.Range("Q18").Value = 610347.1442693399
$ws.Range("R18").Value = 6897598.013066654
# Row 19
$ws.Range("A19").Value = 111486580
$ws.Range("B19").Value = 96348
$ws.Range("D19").Value = "VU"
$ws.Range("E19").Value = 220787
$ws.Range("F19").Value = "Knärot"
$ws.Range("G19").Value = "Goodyera repens"
$ws.Range("H19").Value = "(L.) R. Br."
$ws.Range("Q19").Value = 610100.4635512675
$ws.Range("R19").Value = 6897186.766084836
# Row 20
$ws.Range("A20").Value = 111486557
$ws.Range("Q20").Value = 610221.9365824561
$ws.Range("R20").Value = 6897535.519929474
# Row 21
$ws.Range("A21").Value = 111486575
$ws.Range("Q21").Value = 610234.356536509
$ws.Range("R21").Value = 6897284.584036393
# Row 22
$ws.Range("A22").Value = 111486547
$ws.Range("B22").Value = 76495
$ws.Range("E22").Value = 6487
$ws.Range("F22").Value = "Blågrå svartspik"
$ws.Range("G22").Value = "Chaenothecopsis fennica"
$ws.Range("H22").Value = "(Laurila) Tibell"
$ws.Range("Q22").Value = 610214.3838179935
$ws.Range("R22").Value = 6897609.083555548
# Row 24
$ws.Range("A24").Value = 111486568
$ws.Range("B24").Value = 93057
$ws.Range("D24").Value = "LC"
$ws.Range("E24").Value = 2809
$ws.Range("F24").Value = "Mörk husmossa"
$ws.Range("G24").Value = "Hylocomiastrum umbratum"
$ws.Range("H24").Value = "(Hedw.) M.Fleisch."
$ws.Range("Q24").Value = 610178.0574054071
$ws.Range("R24").Value = 6897403.74427297
# Row 25
$ws.Range("A25").Value = 111486560
$ws.Range("B25").Value = 89405
$ws.Range("D25").Value = "NT"
$ws.Range("E25").Value = 1202
$ws.Range("F25").Value = "Ullticka"
$ws.Range("G25").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H25").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q25").Value = 610242.9519009186
$ws.Range("R25").Value = 6897536.674650873
$ws.Range("AC25").Value = ""
# Row 26
$ws.Range("A26").Value = 111486548
$ws.Range("B26").Value = 77515
$ws.Range("E26").Value = 6425
$ws.Range("F26").Value = "Garnlav"
$ws.Range("G26").Value = "Alectoria sarmentosa"
$ws.Range("H26").Value = "(Ach.) Ach."
$ws.Range("Q26").Value = 610211.0056098022
$ws.Range("R26").Value = 6897612.243104065
# Row 27
$ws.Range("A27").Value = 111486583
$ws.Range("B27").Value = 89405
$ws.Range("E27").Value = 1202
$ws.Range("F27").Value = "Ullticka"
$ws.Range("G27").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H27").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q27").Value = 610170.0873515971
$ws.Range("R27").Value = 6897389.935445569
# Row 28
$ws.Range("A28").Value = 111486554
$ws.Range("B28").Value = 78107
$ws.Range("E28").Value = 6453
$ws.Range("F28").Value = "Vedskivlav"
$ws.Range("G28").Value = "Hertelidea botryosa"
$ws.Range("H28").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q28").Value = 610206.4050188576
$ws.Range("R28").Value = 6897581.260525526
# Row 29
$ws.Range("A29").Value = 111486551
$ws.Range("B29").Value = 89425
$ws.Range("D29").Value = "NT"
$ws.Range("E29").Value = 5442
$ws.Range("F29").Value = "Tallticka"
$ws.Range("G29").Value = "Porodaedalea pini"
$ws.Range("H29").Value = "(Brot.) Murrill"
$ws.Range("Q29").Value = 610206.9311065455
$ws.Range("R29").Value = 6897622.387123355
# Row 30
$ws.Range("A30").Value = 111486556
$ws.Range("B30").Value = 78107
$ws.Range("D30").Value = "NT"
$ws.Range("E30").Value = 6453
$ws.Range("F30").Value = "Vedskivlav"
$ws.Range("G30").Value = "Hertelidea botryosa"
$ws.Range("H30").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q30").Value = 610193.6094834032
$ws.Range("R30").Value = 6897571.966054032
# Row 31
$ws.Range("A31").Value = 111486574
$ws.Range("B31").Value = 96348
$ws.Range("D31").Value = "VU"
$ws.Range("E31").Value = 220787
$ws.Range("F31").Value = "Knärot"
$ws.Range("G31").Value = "Goodyera repens"
$ws.Range("H31").Value = "(L.) R. Br."
$ws.Range("Q31").Value = 610224.6307519674
$ws.Range("R31").Value = 6897324.444504307
# Row 32
$ws.Range("A32").Value = 111486576
$ws.Range("B32").Value = 77186
$ws.Range("D32").Value = "NT"
$ws.Range("E32").Value = 353
$ws.Range("F32").Value = "Dvärgbägarlav"
$ws.Range("G32").Value = "Cladonia parasitica"
$ws.Range("H32").Value = "(Hoffm.) Hoffm."
$ws.Range("Q32").Value = 610214.2438761768
$ws.Range("R32").Value = 6897284.393316317
# Row 33
$ws.Range("A33").Value = 111486581
$ws.Range("Q33").Value = 610072.3796948178
$ws.Range("R33").Value = 6897129.783162965
# Row 34
$ws.Range("A34").Value = 111486584
$ws.Range("B34").Value = 78512
$ws.Range("D34").Value = "LC"
$ws.Range("E34").Value = 6456
$ws.Range("F34").Value = "Skinnlav"
$ws.Range("G34").Value = "Leptogium saturninum"
$ws.Range("H34").Value = "(Dicks.) Nyl."
$ws.Range("Q34").Value = 610335.6047682473
$ws.Range("R34").Value = 6897578.948932082

Write-Output "Applied row reshuffle changes to rows 2-34 (except 6 and 23)"
